$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value updates (no row movement) ---

# Row 9: Alemania
$ws.Range("E9").Value = 33749
$ws.Range("G9").Value = 16
$ws.Range("H9").Value = 6330

# Row 17: Paises Bajos
$ws.Range("F17").Value = 804

# Row 121: Estado de Palestina
$ws.Range("B121").Value = 344
$ws.Range("C121").Value = 1
$ws.Range("E121").Value = 271

# Row 122: Mauricio
$ws.Range("D122").Value = 306
$ws.Range("E122").Value = 18

# --- Row re-ranking: countries move up due to case-count increase ---
# Pattern: delete the country's old row, then insert a fresh row at its
# new (higher) position and fill it with the updated data. This causes the
# rows in between to shift down by one, matching the diff exactly.

# Moldavia: was row 60 -> now row 57 (between Argentina and Luxemburgo)
$ws.Rows(60).Delete()
$ws.Rows(57).Insert()
$ws.Range("A57").Value = "Moldavia"
$ws.Range("B57").Value = 3771
$ws.Range("C57").Value = 133
$ws.Range("D57").Value = 1114
$ws.Range("E57").Value = 2550
$ws.Range("F57").Value = 212
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 107

# Somalia: was row 108 -> now row 106 (between Sri Lanka and San Marino)
$ws.Rows(108).Delete()
$ws.Rows(106).Insert()
$ws.Range("A106").Value = "Somalia"
$ws.Range("B106").Value = 582
$ws.Range("C106").Value = 54
$ws.Range("D106").Value = 20
$ws.Range("E106").Value = 534
$ws.Range("F106").Value = 2
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 28

# Togo: was row 150 -> now row 148 (between Bermudas and Sierra Leona)
$ws.Rows(150).Delete()
$ws.Rows(148).Insert()
$ws.Range("A148").Value = "Togo"
$ws.Range("B148").Value = 109
$ws.Range("C148").Value = 10
$ws.Range("D148").Value = 64
$ws.Range("E148").Value = 38
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 1
$ws.Range("H148").Value = 7

# --- Update the "last updated" timestamp text in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 16:52"

Write-Host "done"
